$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells are stored as text in this sheet (values such as
# "264.30" / "0.0470" / "51.208.11" must keep their exact printed form, not be
# coerced to numbers), so a leading apostrophe forces text entry for each one.

$ws.Range("D2").Value = "'51.208.11"

$ws.Range("D3").Value = "'3.064.52"
$ws.Range("E3").Value = "  +1.32%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'387.89"
$ws.Range("E5").Value = "  +2.14%  "

$ws.Range("D6").Value = "'102.12"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("E7").Value = "  -1.76%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "'0.579"
$ws.Range("E9").Value = "  -1.61%  "

$ws.Range("E10").Value = "  +0.30%  "

$ws.Range("E11").Value = "  +0.34%  "

$ws.Range("E12").Value = "  -1.24%  "

$ws.Range("D13").Value = "'3.545.53"
$ws.Range("E13").Value = "  +1.19%  "

$ws.Range("D14").Value = "'18.27"
$ws.Range("E14").Value = "  -1.33%  "

$ws.Range("D15").Value = "'7.69"
$ws.Range("E15").Value = "  -0.51%  "

$ws.Range("D16").Value = "'3.070.20"
$ws.Range("E16").Value = "  +1.69%  "

$ws.Range("D17").Value = "'0.995"
$ws.Range("E17").Value = "  +2.24%  "

$ws.Range("D18").Value = "'10.71"
$ws.Range("E18").Value = "  +0.93%  "

$ws.Range("D19").Value = "'51.219.57"
$ws.Range("E19").Value = "  -0.76%  "

$ws.Range("E20").Value = "  +2.69%  "

$ws.Range("D21").Value = "'12.26"
$ws.Range("E21").Value = "  -1.27%  "

$ws.Range("D22").Value = "'0.0₃0955"
$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("E23").Value = "  -0.35%  "

$ws.Range("D24").Value = "'264.30"
$ws.Range("E24").Value = "  -1.09%  "

$ws.Range("E25").Value = "  -1.10%  "

$ws.Range("D26").Value = "'7.89"
$ws.Range("E26").Value = "  -6.81%  "

$ws.Range("D27").Value = "'26.87"
$ws.Range("E27").Value = "  +2.76%  "

$ws.Range("D28").Value = "'7.29"
$ws.Range("E28").Value = "  -2.83%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("E30").Value = "  -5.26%  "

$ws.Range("E31").Value = "  -3.19%  "

$ws.Range("E32").Value = "  +1.43%  "

$ws.Range("D33").Value = "'35.58"
$ws.Range("E33").Value = "  +4.63%  "

$ws.Range("D34").Value = "'0.0470"
$ws.Range("E34").Value = "  +4.71%  "

$ws.Range("D35").Value = "'2.07"
$ws.Range("E35").Value = "  +2.50%  "

$ws.Range("D36").Value = "'49.99"
$ws.Range("E36").Value = "  -1.17%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").Value = "'3.35"
$ws.Range("E38").Value = "  +1.17%  "

$ws.Range("D39").Value = "'0.295"
$ws.Range("E39").Value = "  -1.56%  "

$ws.Range("D40").Value = "'130.49"
$ws.Range("E40").Value = "  +1.18%  "

$ws.Range("E41").Value = "  -3.27%  "

$ws.Range("E42").Value = "  -1.59%  "

$ws.Range("E43").Value = "  -0.96%  "

$ws.Range("D44").Value = "'2.48"
$ws.Range("E44").Value = "  -2.16%  "

$ws.Range("D45").Value = "'3.71"
$ws.Range("E45").Value = "  -0.79%  "

$ws.Range("D46").Value = "'21.69"

$ws.Range("E47").Value = "  +3.52%  "

$ws.Range("E48").Value = "  -0.18%  "

$ws.Range("D49").Value = "'2.066.63"
$ws.Range("E49").Value = "  +2.30%  "

$ws.Range("D50").Value = "'0.0326"
$ws.Range("E50").Value = "  +3.91%  "

$ws.Range("D51").Value = "'0.903"
$ws.Range("E51").Value = "  +14.22%  "
